# Daily attendance processing - 2025-11-28 09:50:37
# Normalize the "Recorded By" (column G) entries: when the list of
# recorders (comma-separated) starts with the literal token "System",
# move it to the end of the list instead (reverse the ordering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Length -gt 1 -and $parts[0].Equals("System")) {
            $n = $parts.Length
            $rev = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $rev += $parts[$i]
            }
            $newVal = $rev -join ", "
            $cell.Value2 = $newVal
        }
    }
}
